$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet) - update "想去人数" (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 86
$ws1.Range("F5").Value = 26
$ws1.Range("F6").Value = 23
$ws1.Range("F7").Value = 563
$ws1.Range("F8").Value = 7979
$ws1.Range("F9").Value = 753
$ws1.Range("F10").Value = 229
$ws1.Range("F11").Value = 1096
$ws1.Range("F12").Value = 762
$ws1.Range("F13").Value = 29
$ws1.Range("F15").Value = 200
$ws1.Range("F16").Value = 26
$ws1.Range("F17").Value = 47
$ws1.Range("F18").Value = 206
$ws1.Range("F19").Value = 830

# Sheet "全部类型" (fourth sheet) - update "想去人数" (F column) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 86
$ws4.Range("F5").Value = 26
$ws4.Range("F6").Value = 23
$ws4.Range("F8").Value = 563
$ws4.Range("F9").Value = 7979
$ws4.Range("F10").Value = 753
$ws4.Range("F11").Value = 229
$ws4.Range("F12").Value = 1096
$ws4.Range("F13").Value = 762
$ws4.Range("F14").Value = 29
$ws4.Range("F16").Value = 200
$ws4.Range("F17").Value = 26
$ws4.Range("F18").Value = 47
$ws4.Range("F19").Value = 206
$ws4.Range("F20").Value = 830

$wb.Save()
